$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecast column BB (54): header date, formatted like the BA (53) header
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(1, 54).Value = 45986

# Carry the BA forecast forward into BB for rows 3-18 (unchanged values)
for ($r = 3; $r -le 18; $r++) {
    $baVal = $ws.Cells.Item($r, 53).Value2
    if ($null -ne $baVal) {
        $ws.Cells.Item($r, 54).Value = $baVal
    }
}

# Rows 19-21 carry new, distinct re-forecast values in column BB
$ws.Cells.Item(19, 54).Value = 0.8976398032236155
$ws.Cells.Item(20, 54).Value = -0.6203510926954925
$ws.Cells.Item(21, 54).Value = -1.016686374377895
